# Fruta / hortaliza, semanal
# A new weekly price observation was added to the "Jengibre" sheet.
# It is inserted as a new row 10 (pushing the existing rows 10-70 down
# to 11-71), matching the published diff where every existing row from
# 10 onward shifts down by one and a brand-new row 10 appears with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 10, shifting rows 10..70
# down to 11..71 (this reproduces the Rows.InsertBefore-alike behaviour
# that the diff implies: every row >=10 keeps its original contents but
# moves down one position).
$ws.Rows.Item(10).Insert()

# Populate the freshly inserted row 10 with the new weekly record.
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(10, 3).Value = "Coquimbo"
$ws.Cells.Item(10, 4).Value = 44847
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = 100114007
$ws.Cells.Item(10, 7).Value = "Jengibre"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 400
$ws.Cells.Item(10, 11).Value = 13500
$ws.Cells.Item(10, 12).Value = 14000
$ws.Cells.Item(10, 13).Value = 13750
$ws.Cells.Item(10, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(10, 15).Value = "Perú"
$ws.Cells.Item(10, 16).Value = 1058
$ws.Cells.Item(10, 17).Value = 13
$ws.Cells.Item(10, 18).Value = "Hortaliza"
